# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the regenerated site data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 17
$ws1.Range("F7").Value  = 2645
$ws1.Range("F12").Value = 937
$ws1.Range("F17").Value = 790
$ws1.Range("F24").Value = 226
$ws1.Range("F25").Value = 316
$ws1.Range("F27").Value = 695
$ws1.Range("F28").Value = 571
$ws1.Range("F29").Value = 5674
$ws1.Range("F30").Value = 5674
$ws1.Range("F31").Value = 498
$ws1.Range("F36").Value = 1639

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 87
$ws2.Range("F9").Value  = 42
$ws2.Range("F18").Value = 210

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 17
$ws4.Range("F3").Value  = 87
$ws4.Range("F5").Value  = 2645
$ws4.Range("F10").Value = 937
$ws4.Range("F18").Value = 790
$ws4.Range("F25").Value = 42
$ws4.Range("F28").Value = 316
$ws4.Range("F30").Value = 571
$ws4.Range("F31").Value = 5674
$ws4.Range("F33").Value = 498
$ws4.Range("F37").Value = 1639
